$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 14.47727272727272
$ws.Range("R2").Value = 1.815485677363773
$ws.Range("S2").Value = 1.979371877230549

$ws.Range("K3").Value = 1.791666666666668
$ws.Range("R3").Value = 1.620655622136059
$ws.Range("S3").Value = 1.747323835194455

$ws.Range("K5").Value = 14.47727272727272
$ws.Range("R5").Value = 1.815485677363773
$ws.Range("S5").Value = 1.979371877230549

$ws.Range("K7").Value = 14.47727272727272
$ws.Range("R7").Value = 1.815485677363773
$ws.Range("S7").Value = 1.979371877230549

$ws.Range("K8").Value = 14.47727272727272
$ws.Range("R8").Value = 1.815485677363773
$ws.Range("S8").Value = 1.979371877230549

$ws.Range("K9").Value = 13.76976495726495
$ws.Range("R9").Value = 1.803394296576035
$ws.Range("S9").Value = 1.964819060413116

$ws.Range("K11").Value = 1.791666666666668
$ws.Range("R11").Value = 1.620655622136059
$ws.Range("S11").Value = 1.747323835194455

$ws.Range("K14").Value = 16.97685185185183
$ws.Range("R14").Value = 1.859533546038736
$ws.Range("S14").Value = 2.032558602498382

$ws.Range("K15").Value = -3.847222222222223
$ws.Range("R15").Value = 1.546865537736907
$ws.Range("S15").Value = 1.660778333536659

$ws.Range("K18").Value = 14.47727272727272
$ws.Range("R18").Value = 1.815485677363773
$ws.Range("S18").Value = 1.979371877230549

$ws.Range("K19").Value = 5.462962962962945
$ws.Range("R19").Value = 1.672603071948262
$ws.Range("S19").Value = 1.808689105403011

$ws.Range("K21").Value = 0.2777777777777778
$ws.Range("R21").Value = 1.600162412993039
$ws.Range("S21").Value = 1.723215189873418

$ws.Range("K22").Value = 1.791666666666668
$ws.Range("R22").Value = 1.620655622136059
$ws.Range("S22").Value = 1.747323835194455

$ws.Range("K23").Value = 13.46442495126706
$ws.Range("R23").Value = 1.798225615362447
$ws.Range("S23").Value = 1.958604378795604

$ws.Range("K24").Value = 19.79629629629628
$ws.Range("R24").Value = 1.911855479578636
$ws.Range("S24").Value = 2.09608909874769

$ws.Range("K28").Value = 21.28240740740739
$ws.Range("R28").Value = 1.940636870984383
$ws.Range("S28").Value = 2.131200751448103

$ws.Range("K29").Value = 1.791666666666668
$ws.Range("R29").Value = 1.620655622136059
$ws.Range("S29").Value = 1.747323835194455

$ws.Range("K31").Value = 12.67039049919483

$ws.Range("K32").Value = 12.67039049919483
$ws.Range("R32").Value = 1.784922174701128
$ws.Range("S32").Value = 1.942625691911729

$ws.Range("K33").Value = 19.65277777777778
$ws.Range("R33").Value = 1.909121107266436
$ws.Range("S33").Value = 2.092759415833974

$ws.Range("K34").Value = 21.28240740740739
$ws.Range("R34").Value = 1.940636870984383
$ws.Range("S34").Value = 2.131200751448103

$ws.Range("K35").Value = 13.0158303464755
$ws.Range("R35").Value = 1.790685487585954
$ws.Range("S35").Value = 1.94954496878686

$ws.Range("K36").Value = 13.0158303464755

$ws.Range("K37").Value = 19.60879629629628
$ws.Range("R37").Value = 1.908284719500103
$ws.Range("S37").Value = 2.091741145739967

$ws.Range("K39").Value = 14.47727272727272
$ws.Range("R39").Value = 1.815485677363773
$ws.Range("S39").Value = 1.979371877230549

$ws.Range("K40").Value = 14.47727272727272
$ws.Range("R40").Value = 1.815485677363773
$ws.Range("S40").Value = 1.979371877230549

$ws.Range("K41").Value = 14.96875
$ws.Range("R41").Value = 1.8239809580482
$ws.Range("S41").Value = 1.989608681354817

$ws.Range("K42").Value = 13.46442495126706
$ws.Range("R42").Value = 1.798225615362447
$ws.Range("S42").Value = 1.958604378795604

$ws.Range("K43").Value = 13.46442495126706

$ws.Range("K45").Value = 14.47727272727272
$ws.Range("R45").Value = 1.815485677363773
$ws.Range("S45").Value = 1.979371877230549

$ws.Range("K47").Value = 13.76976495726495
$ws.Range("R47").Value = 1.803394296576035
$ws.Range("S47").Value = 1.964819060413116

$ws.Range("K51").Value = 13.0158303464755
$ws.Range("R51").Value = 1.790685487585954
$ws.Range("S51").Value = 1.94954496878686

$ws.Range("K52").Value = 13.76976495726495
$ws.Range("R52").Value = 1.803394296576035
$ws.Range("S52").Value = 1.964819060413116

$ws.Range("K57").Value = 5.462962962962945
$ws.Range("R57").Value = 1.672603071948262
$ws.Range("S57").Value = 1.979371877230549
